# Remove the "SOBRE" (description) rows from the data dictionary:
# - DESENVOLVEDORAS table: remove the DVD_SOBRE column row
# - DISTRIBUIDORAS table: remove the DBR_SOBRE column row
# These rows are duplicated on the "FULL" summary sheet as well as on the
# individual per-table sheets, so they need to be removed from all three
# places.

$wb = $excel.ActiveWorkbook

# --- FULL sheet: remove DVD_SOBRE row (row 36) and DBR_SOBRE row (row 42) ---
$full = $wb.Sheets.Item("FULL")
# Delete the lower row first so the row index of the first deletion doesn't shift.
$full.Rows.Item(42).Delete() | Out-Null
$full.Rows.Item(36).Delete() | Out-Null

# --- DESENVOLVEDORAS sheet: remove DVD_SOBRE row (row 7) ---
$dvd = $wb.Sheets.Item("DESENVOLVEDORAS")
$dvd.Rows.Item(7).Delete() | Out-Null

# --- DISTRIBUIDORAS sheet: remove DBR_SOBRE row (row 6) ---
$dbr = $wb.Sheets.Item("DISTRIBUIDORAS")
$dbr.Rows.Item(6).Delete() | Out-Null

# --- Restore view/selection state ---
$full.Activate()
$full.Range("A78").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 55

$dvd2 = $wb.Sheets.Item("DESENVOLVEDORAS")
$dvd2.Activate()
$dvd2.Range("E15").Select() | Out-Null

$dbr2 = $wb.Sheets.Item("DISTRIBUIDORAS")
$dbr2.Activate()
$dbr2.Range("C15").Select() | Out-Null

$editoras = $wb.Sheets.Item("EDITORAS")
$editoras.Activate()
$editoras.Range("E8").Select() | Out-Null
